# Applies the "full data scraped for extra batting and bowling fields" update
# to the "ODI Batting Extra" worksheet: 15 additional match rows are inserted
# at the top of the data block (match codes 4064-4457), shifting the existing
# rows down, so the table grows from A1:F21 to A1:F36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

$rows = @(
    @{ Row=2; A="4064"; B=11; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=3; A="4065"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=4; A="4135"; B=9; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=5; A="4175"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=6; A="4196"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=7; A="4387"; B=10; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=8; A="4388"; B=8; C="2"; D="1"; E="5.38%"; F="NO" },
    @{ Row=9; A="4398"; B=8; C="2"; D="0"; E="5.10%"; F="NO" },
    @{ Row=10; A="4402"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=11; A="4406"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=12; A="4410"; B=8; C="1"; D="0"; E="2.36%"; F="NO" },
    @{ Row=13; A="4437"; B=8; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=14; A="4454"; B=8; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=15; A="4456"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=16; A="4457"; B=8; C="1"; D="3"; E="9.12%"; F="NO" },
    @{ Row=17; A="4524"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=18; A="4526"; B=7; C="3"; D="1"; E="13.94%"; F="NO" },
    @{ Row=19; A="4533"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=20; A="4535"; B=8; C="1"; D="0"; E="3.38%"; F="NO" },
    @{ Row=21; A="4621"; B=8; C="1"; D="0"; E="2.27%"; F="NO" },
    @{ Row=22; A="4623"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=23; A="4624"; B=8; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=24; A="4640"; B=8; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=25; A="4643"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=26; A="4656"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=27; A="4657"; B=6; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=28; A="4658"; B=6; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=29; A="4669"; B=8; C="0"; D="0"; E="0.33%"; F="NO" },
    @{ Row=30; A="4679"; B=8; C="0"; D="0"; E="1.08%"; F="NO" },
    @{ Row=31; A="4682"; B=7; C="0"; D="0"; E="2.63%"; F="NO" },
    @{ Row=32; A="4685"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=33; A="4692"; B=$null; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=34; A="4695"; B=8; C=$null; D=$null; E=$null; F="NO" },
    @{ Row=35; A="4697"; B=8; C="3"; D="1"; E="6.49%"; F="YES" },
    @{ Row=36; A="4725"; B=8; C=$null; D=$null; E=$null; F="NO" }
)

# Columns A, C, D and E hold values that look numeric ("4064", "0", "1",
# "13.94%", ...) but must stay stored as text. Pre-format the whole block as
# text so assigning the values does not get auto-converted into numbers or
# percentages by Excel. Column B (BATTING_POSITION) holds real numbers and is
# left with its default (General) number format.
$ws.Range("A2:A36").NumberFormat = "@"
$ws.Range("C2:E36").NumberFormat = "@"

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
